# Added proper comment in Write_data_to_csc file
# Append two new employee records to the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Employee number" column holds values that look numeric
# (e.g. "543262") but must be stored as text, matching the rest
# of the column. Mark the range as Text before typing the values
# so Excel keeps them as strings instead of coercing to numbers,
# then drop the formatting again so the cells fall back to the
# sheet's default (unstyled) appearance.
$numberRange = $ws.Range("B11:B12")
$numberRange.NumberFormat = "@"

$ws.Range("A11").Value = "salmanalam"
$ws.Range("B11").Value = "543262"
$ws.Range("C11").Value = "Devops"
$ws.Range("D11").Value = "m"

$ws.Range("A12").Value = "farmankhan"
$ws.Range("B12").Value = "345672"
$ws.Range("C12").Value = "Ba"
$ws.Range("D12").Value = "m"

$numberRange.ClearFormats()
